$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Strategy 2 / Unigram Strict DF): mark preprocessing step as "Model trained"
# and note that it was sent to Sam
$ws.Range("F4").Value = "2. Model trained"
$ws.Range("G4").Value = "(sent to Sam)"

# Row 6 (Strategy 4 / Bigram Base): already "Model trained" - add note it was sent to Sam
$ws.Range("G6").Value = "(sent to Sam)"

# Row 7 (Strategy 5 / Bigram Strict DF): already "Model trained" - add note it was sent to Sam
$ws.Range("G7").Value = "(sent to Sam)"

# Row 8 (Strategy 6 / Bigram Proper Only): mark as processed/trained, and update the note
# to reflect that it still needs to be sent to Sam
$ws.Range("F8").Value = "2. Model trained"
$ws.Range("G8").Value = "(need to send to Sam)"

# F8 now has a value, so drop it from the dropdown data validation range (F3 F5:F6 F8 -> F3 F5:F6)
$ws.Range("F8").Validation.Delete()

# Update the active selection to reflect where editing left off
$ws.Range("G8").Select()
